$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Purchase 22-23")
$ws2 = $wb.Worksheets.Item("Sale 22-23")

# --- Add the new purchase row (Sr. No 6 / Cassun Electricals) to "Purchase 22-23" ---
# Duplicate the formatting of the row above (Sr. No 4, row 13) down onto row 17,
# then overwrite the copied values with the new record's data.
$ws1.Range("A13:F13").Copy($ws1.Range("A17:F17")) | Out-Null

$ws1.Range("A17").Value = 6
$ws1.Range("B17").Value = 45234
$ws1.Range("C17").Value = "INV/23-24/1281"
$ws1.Range("D17").Value = "Cassun Electricals"
$ws1.Range("E17").Value = 14313
$ws1.Range("F17").Formula = "=E17"

# --- Update sheet selections / active tab ---
# "Sale 22-23" ends up with a plain (non-active) selection over A45:E49.
$ws2.Range("A45:E49").Select() | Out-Null

# "Purchase 22-23" becomes the active / tab-selected sheet, with B23 selected.
$ws1.Activate() | Out-Null
$ws1.Range("B23").Select() | Out-Null
